# Atualizando base de dados da pesquisa via Streamlit
# Adds a new survey response row (row 16) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 16

# Columns that stay blank in this response still need a present-but-empty
# cell (matching how the rest of the sheet was exported), so nudge a
# (no-op) formatting property on each to force Excel to materialize it.
$blankCols = 2,3,4,5,6,7,8,10,11,12,13,14,15,16,17,18,19,20
foreach ($col in $blankCols) {
    $cell = $ws.Cells.Item($newRow, $col)
    $cell.Value = ""
    $cell.Font.Bold = $false
}

$ws.Cells.Item($newRow, 1).Value = "jessica.mariano@mrv.com.br"
$ws.Cells.Item($newRow, 9).Value = "2025-05-20 20:33:52"
$ws.Cells.Item($newRow, 21).Value = "Painel do Portifólio - Planejamento da Produção - PLNESROBR004: Trazer visão de unidades na aba PEI também."
$ws.Cells.Item($newRow, 22).Value = "Planilha de Médio Prazo (Replan),Acompanhamento das metas mensais do engenheiro,Excel,💎 Muito Importante,40.0"
